$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.794.55'
$ws.Range("E2").Value = '  -2.47%  '
$ws.Range("D3").Value = '3.273.73'
$ws.Range("E3").Value = '  -1.20%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.17'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.18%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  -3.02%  '
$ws.Range("E10").Value = '  +0.43%  '
$ws.Range("E11").Value = '  -2.75%  '
$ws.Range("D12").Value = '3.845.95'
$ws.Range("E12").Value = '  -1.08%  '
$ws.Range("E13").Value = '  -3.84%  '
$ws.Range("D14").Value = '65.874.58'
$ws.Range("E14").Value = '  -2.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.42'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.74%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000162'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.89%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.264.56'
$ws.Range("E17").Value = '  -1.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '436.09'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.83%  '
$ws.Range("E19").Value = '  -2.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.38'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.34'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").Value = '3.424.90'
$ws.Range("E24").Value = '  -0.90%  '
$ws.Range("E25").Value = '  -2.14%  '
$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000113'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.15%  '
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.195'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.87'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.88%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("E30").Value = '  -2.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.23'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.01%  '
$ws.Range("E32").Value = '  +0.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.13'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.85%  '
$ws.Range("E34").Value = '  -3.24%  '
$ws.Range("E35").Value = '  -5.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '159.04'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.41%  '
$ws.Range("E37").Value = '  -5.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '26.73'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.78'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.94%  '
$ws.Range("D40").Value = '2.772.12'
$ws.Range("E40").Value = '  +0.49%  '
$ws.Range("E41").Value = '  -2.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.31'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.27'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.02'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0655'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.66%  '
$ws.Range("E46").Value = '  -5.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '317.91'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.32'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.05%  '
$ws.Range("E49").Value = '  -2.28%  '
$ws.Range("E50").Value = '  +2.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.04%  '
